$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the entire contents of row 3 and row 4 (all columns),
# including the optional "M" (Aktivitet) column which only row 4 had.

$cols = @("A","B","C","D","E","F","G","H","M","P","Q","R","S")

foreach ($col in $cols) {
    $addr3 = $col + "3"
    $addr4 = $col + "4"
    $v3 = $ws.Range($addr3).Value2
    $v4 = $ws.Range($addr4).Value2
    $ws.Range($addr3).Value = $v4
    $ws.Range($addr4).Value = $v3
}

# Row 4 no longer has a value in column M after the swap (row 3's original
# M cell was empty), so clear it out.
$ws.Range("M4").ClearContents()
